# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G ("K" - strikeouts) is recalculated for each data row (rows 2-65)
# on Sheet1, replacing the previous values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G, rows 2 through 65 (in row order)
$newK = @(
    0,0,2,1,1,1,0,0,2,1,
    2,1,2,4,2,1,1,2,0,1,
    2,0,0,3,1,1,1,1,0,1,
    2,3,1,1,0,1,2,1,4,2,
    0,1,0,0,2,1,1,2,1,0,
    0,2,2,2,1,5,2,4,5,1,
    2,0,2,1
)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
